# prometheus-net/Docs/MeasurementsBenchmarks.xlsx
# Update the four raw latency measurements (column G) on Sheet1 with a
# fresh benchmark run, then leave the active selection where the next
# entry would go (one row below the data, over in column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 383.5
$ws.Range("G4").Value = 169.1
$ws.Range("G5").Value = 953.5
$ws.Range("G6").Value = 1531.1

$ws.Range("G7").Select()
